$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.0
$ws.Range("A7").Value = 9.0
$ws.Range("A8").Value = 10.500000000000002
$ws.Range("A11").Value = 9.0
$ws.Range("A13").Value = 0.0
$ws.Range("A16").Value = 3.0
$ws.Range("A18").Value = 0.0
$ws.Range("A22").Value = 6.0
